$d = $word.ActiveDocument

# 1) Footnote "1. ^ 萌铺子..." paragraph mark: add an eastAsia font hint.
#    (No visible text change - this only touches paragraph mark run formatting,
#    so we find the paragraph by its body text and flip the paragraph-mark font.)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*萌铺子（杭州）科技有限公司 遵循 先立后破 不立不破 原则*") {
        $p.Range.Font.NameFarEast = "eastAsia-hint"
        break
    }
}

# 2) Footnote 6 paragraph: "...股东代表大会 / 职工代表大会 [6]" -> "...[7]"
$d.Content.Find.Execute(
    "企业 会议 遵循 全过程人民民主 原则 有 股东代表大会 / 职工代表大会 [6]", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "企业 会议 遵循 全过程人民民主 原则 有 股东代表大会 / 职工代表大会 [7]", 2)

# 3) Footnote 7 paragraph: "...负责 [7]" -> "...负责 [8]"
$d.Content.Find.Execute(
    "股东代表大会 职责 是 接收 职工代表大会 职工代表 (高级管理人员) 报告 / 由 股东代表大会 股东代表 (法定代表人) 依 法定 程序 办理 事项 且 记录 / 存档 / 修订 / 审议 章程 因 股东代表大会 股东 表决权 产生 对 股东代表大会 股东 表决权 负责 股东代表大会 因 中华人民共和国公司法 产生 对 中华人民共和国公司法 负责 [7]", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "股东代表大会 职责 是 接收 职工代表大会 职工代表 (高级管理人员) 报告 / 由 股东代表大会 股东代表 (法定代表人) 依 法定 程序 办理 事项 且 记录 / 存档 / 修订 / 审议 章程 因 股东代表大会 股东 表决权 产生 对 股东代表大会 股东 表决权 负责 股东代表大会 因 中华人民共和国公司法 产生 对 中华人民共和国公司法 负责 [8]", 2)

# 4) Footnote 11 paragraph (English): "...startups [9]" -> "...startups [12]"
$d.Content.Find.Execute(
    "Serving legal representatives of startups in their ultra-early phase / Services for legal representatives of ultra-early-phase startups [9]", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Serving legal representatives of startups in their ultra-early phase / Services for legal representatives of ultra-early-phase startups [12]", 2)

# 5) Header date field cached text update.
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute(
        "2025/05/26 12:39 PM", $true, $false, $false, $false, $false,
        $true, 1, $false,
        "2025/05/29 01:39 PM", 2)
}
